$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-19 Tuesday" "2023-09-20 Wednesday"

Replace-Text "89×25=2225" "79×42=3318"
Replace-Text "13×12=156" "81×26=2106"
Replace-Text "34×16=544" "50×37=1850"
Replace-Text "89×73=6497" "83×56=4648"
Replace-Text "68×22=1496" "73×68=4964"

Replace-Text "21×92=1932" "73×58=4234"
Replace-Text "68×31=2108" "26×43=1118"
Replace-Text "22×20=440" "29×22=638"
Replace-Text "39×26=1014" "59×45=2655"
Replace-Text "15×38=570" "40×38=1520"

Replace-Text "80×49=3920" "87×32=2784"
Replace-Text "20×86=1720" "66×91=6006"
Replace-Text "80×48=3840" "29×42=1218"
Replace-Text "66×31=2046" "13×30=390"
Replace-Text "43×72=3096" "21×57=1197"

Replace-Text "81×46=3726" "89×27=2403"
Replace-Text "87×14=1218" "95×38=3610"
Replace-Text "47×91=4277" "85×32=2720"
Replace-Text "84×90=7560" "21×68=1428"
Replace-Text "72×76=5472" "59×99=5841"

Replace-Text "46×56=2576" "24×67=1608"
Replace-Text "21×55=1155" "31×82=2542"
Replace-Text "26×12=312" "15×61=915"
Replace-Text "25×54=1350" "13×69=897"
Replace-Text "24×73=1752" "25×31=775"
